# Scheduled runner: refresh cached Market Board pricing columns (H:N)
# across the Leve-profit sheets. Values come from the latest Universalis
# snapshot; only currentAveragePrice* / LevePrice* / LeveProfit* cells move.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 51.666668
$ws.Range("I8").Value = 51.666668
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 155.000004
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -16.00000399999999
$ws.Range("N8").ClearContents()

# Row 18
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -716
$ws.Range("N18").ClearContents()

# Row 39
$ws.Range("H39").Value = 149.83333
$ws.Range("I39").Value = 94.5
$ws.Range("J39").Value = 194.1
$ws.Range("K39").Value = 283.5
$ws.Range("L39").Value = 582.3
$ws.Range("M39").Value = 12.5
$ws.Range("N39").Value = -1174.3

# Row 40
$ws.Range("H40").Value = 1453.3334
$ws.Range("I40").Value = 1100.125
$ws.Range("J40").Value = 1857
$ws.Range("K40").Value = 1100.125
$ws.Range("L40").Value = 1857
$ws.Range("M40").Value = -925.125
$ws.Range("N40").Value = -2207

# Row 53
$ws.Range("H53").Value = 3338.2
$ws.Range("I53").Value = 350
$ws.Range("J53").Value = 5330.3335
$ws.Range("K53").Value = 350
$ws.Range("L53").Value = 5330.3335
$ws.Range("M53").Value = 287
$ws.Range("N53").Value = -6604.3335

# Row 61
$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

# Row 112
$ws.Range("H112").Value = 3832484.5
$ws.Range("I112").Value = 599.5
$ws.Range("J112").Value = 4116327.8
$ws.Range("K112").Value = 1798.5
$ws.Range("L112").Value = 12348983.4
$ws.Range("M112").Value = -690.5
$ws.Range("N112").Value = -12351199.4

# Row 116
$ws.Range("H116").Value = 4326.3687
$ws.Range("I116").Value = 2730.7778
$ws.Range("K116").Value = 2730.7778
$ws.Range("M116").Value = 711.2222000000002

# Row 125
$ws.Range("H125").Value = 2295.375
$ws.Range("I125").Value = 1758
$ws.Range("J125").Value = 2832.75
$ws.Range("K125").Value = 15822
$ws.Range("L125").Value = 25494.75
$ws.Range("M125").Value = -13362
$ws.Range("N125").Value = -30414.75

# Row 129
$ws.Range("H129").Value = 1172.8395
$ws.Range("J129").Value = 1203.2987
$ws.Range("L129").Value = 3609.8961
$ws.Range("N129").Value = -13609.8961

# Row 132
$ws.Range("H132").Value = 3777.2856
$ws.Range("I132").Value = 3701.2104
$ws.Range("K132").Value = 11103.6312
$ws.Range("M132").Value = -8573.6312

# Row 138
$ws.Range("H138").Value = 1298.4865
$ws.Range("I138").Value = 495.6154
$ws.Range("J138").Value = 3196.182
$ws.Range("K138").Value = 1486.8462
$ws.Range("L138").Value = 9588.545999999998
$ws.Range("M138").Value = 3653.1538
$ws.Range("N138").Value = -19868.546

# Row 141
$ws.Range("H141").Value = 3217.2727
$ws.Range("I141").Value = 2378
$ws.Range("J141").Value = 3916.6667
$ws.Range("K141").Value = 7134
$ws.Range("L141").Value = 11750.0001
$ws.Range("M141").Value = -1954
$ws.Range("N141").Value = -22110.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 70081.94
$ws.Range("I32").Value = 72419.8
$ws.Range("K32").Value = 72419.8
$ws.Range("M32").Value = -72132.8

# Row 110
$ws.Range("H110").Value = 498.25
$ws.Range("I110").Value = 444.33334
$ws.Range("K110").Value = 444.33334
$ws.Range("M110").Value = 1600.66666

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2612.7778
$ws.Range("I20").Value = 2394.1667
$ws.Range("J20").Value = 3050
$ws.Range("K20").Value = 2394.1667
$ws.Range("L20").Value = 3050
$ws.Range("M20").Value = -2147.1667
$ws.Range("N20").Value = -3544

# Row 35
$ws.Range("H35").Value = 29800
$ws.Range("J35").Value = 29800
$ws.Range("L35").Value = 29800
$ws.Range("N35").Value = -30420

# Row 94
$ws.Range("H94").Value = 1339.907
$ws.Range("I94").Value = 976.35297
$ws.Range("J94").Value = 2713.3333
$ws.Range("K94").Value = 976.35297
$ws.Range("L94").Value = 2713.3333
$ws.Range("M94").Value = -525.35297
$ws.Range("N94").Value = -3615.3333

# Row 105
$ws.Range("H105").Value = 3575430.8
$ws.Range("I105").Value = 5118.5713
$ws.Range("J105").Value = 7145743
$ws.Range("K105").Value = 5118.5713
$ws.Range("L105").Value = 7145743
$ws.Range("M105").Value = -3371.5713
$ws.Range("N105").Value = -7149237

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 22838.174
$ws.Range("I58").Value = 1265.4445
$ws.Range("J58").Value = 100500
$ws.Range("K58").Value = 1265.4445
$ws.Range("L58").Value = 100500
$ws.Range("M58").Value = -1062.4445
$ws.Range("N58").Value = -100906

# Row 105
$ws.Range("H105").Value = 7354115.5
$ws.Range("I105").Value = 17858052
$ws.Range("J105").Value = 1359.9
$ws.Range("K105").Value = 17858052
$ws.Range("L105").Value = 1359.9
$ws.Range("M105").Value = -17856305
$ws.Range("N105").Value = -4853.9

# Row 132
$ws.Range("H132").Value = 27603.8
$ws.Range("I132").Value = 42905.5
$ws.Range("J132").Value = 4651.25
$ws.Range("K132").Value = 128716.5
$ws.Range("L132").Value = 13953.75
$ws.Range("M132").Value = -126186.5
$ws.Range("N132").Value = -19013.75

# Row 136
$ws.Range("H136").Value = 22838.174
$ws.Range("I136").Value = 1265.4445
$ws.Range("J136").Value = 100500
$ws.Range("K136").Value = 3796.3335
$ws.Range("L136").Value = 301500
$ws.Range("M136").Value = -1246.3335
$ws.Range("N136").Value = -306600

$ws = $wb.Worksheets.Item("CUL")
# Row 26
$ws.Range("H26").Value = 459
$ws.Range("I26").Value = 153.75
$ws.Range("J26").Value = 703.2
$ws.Range("K26").Value = 461.25
$ws.Range("L26").Value = 2109.6
$ws.Range("M26").Value = -173.25
$ws.Range("N26").Value = -2685.6

# Row 131
$ws.Range("H131").Value = 775
$ws.Range("I131").Value = 580.8333
$ws.Range("J131").Value = 787.3936
$ws.Range("K131").Value = 1742.4999
$ws.Range("L131").Value = 2362.1808
$ws.Range("M131").Value = 3297.5001
$ws.Range("N131").Value = -12442.1808

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 65399.457
$ws.Range("I132").Value = 52978.7
$ws.Range("K132").Value = 158936.1
$ws.Range("M132").Value = -156406.1

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 2171.6365
$ws.Range("I93").Value = 2288.8
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 2288.8
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -1040.8
$ws.Range("N93").Value = -3496

# Row 100
$ws.Range("H100").Value = 2011.6471
$ws.Range("I100").Value = 1810.3
$ws.Range("J100").Value = 2299.2856
$ws.Range("K100").Value = 1810.3
$ws.Range("L100").Value = 2299.2856
$ws.Range("M100").Value = -1269.3
$ws.Range("N100").Value = -3381.2856

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1690263.5
$ws.Range("J113").Value = 4505021.5
$ws.Range("L113").Value = 13515064.5
$ws.Range("N113").Value = -13519404.5
